$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text explicitly (NumberFormat "@") so Excel does not auto-convert them
# to numeric values, matching the original inlineStr/text cell content.
# The style is reset back to "Normal" afterwards so no stray formatting
# is left behind on the cell.
$textCells = "D5", "D7", "D10", "D12", "D13", "D14", "D16", "D21", "D22", "D26", "D28", "D29", "D30", "D32", "D35", "D39", "D40", "D44", "D47", "D48", "D49", "D50", "D51"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "36.662.73"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.963.40"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "244.63"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("D7").Value = "60.77"
$ws.Range("E7").Value = "  +8.07%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").Value = "0.0798"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "14.31"
$ws.Range("E12").Value = "  +7.69%  "
$ws.Range("D13").Value = "0.840"
$ws.Range("E13").Value = "  +5.26%  "
$ws.Range("D14").Value = "21.87"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").Value = "2.251.46"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").Value = "1.963.08"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "36.576.65"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "230.51"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "5.11"
$ws.Range("E22").Value = "  +3.95%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +6.63%  "
$ws.Range("E25").Value = "  +4.07%  "
$ws.Range("D26").Value = "0.144"
$ws.Range("E26").Value = "  +10.36%  "
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("D28").Value = "160.86"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "19.44"
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  +12.37%  "
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").Value = "4.78"
$ws.Range("E32").Value = "  +6.25%  "
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  +8.09%  "
$ws.Range("D35").Value = "3.61"
$ws.Range("E35").Value = "  +22.09%  "
$ws.Range("E36").Value = "  +6.78%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").Value = "5.59"
$ws.Range("E39").Value = "  -5.41%  "
$ws.Range("D40").Value = "0.0988"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("E42").Value = "  +3.57%  "
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("D44").Value = "16.34"
$ws.Range("E44").Value = "  +5.31%  "
$ws.Range("D45").Value = "1.368.59"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("D47").Value = "88.76"
$ws.Range("E47").Value = "  +4.52%  "
$ws.Range("D48").Value = "7.18"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").Value = "2.85"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "44.46"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").Value = "1.85"
$ws.Range("E51").Value = "  +6.95%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
